# Applies the Dec-10-2023 cryptos list refresh (prices + 1h volumes),
# including the NEARProtocol / TrustWalletToken row swap at rows 47-48.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.897.38"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").Value = "2.353.03"
$ws.Range("E3").Value = "  -0.57%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  -2.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "240.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "73.44"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.602"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").Value = "  -2.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "59.21"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.43%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "33.53"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +5.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "2.703.51"
$ws.Range("E15").Value = "  -0.54%  "
$ws.Range("E16").Value = "  -3.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.913"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.60%  "
$ws.Range("D18").Value = "2.354.74"
$ws.Range("E18").Value = "  -0.13%  "
$ws.Range("D19").Value = "43.784.36"
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.73"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "77.58"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.30%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "257.21"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +17.54%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.75"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.51"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.59"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.21%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.23"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.60%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.20%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "177.85"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.48%  "
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.136"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0759"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.21"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.45"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("E37").Value = "  -1.86%  "
$ws.Range("E38").Value = "  -3.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.57%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "67.11"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +25.80%  "
$ws.Range("E42").Value = "  +10.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.74%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.202"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.93%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.51"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  -1.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.18%  "
